$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Move the two trailing "meta" columns (old K:L -> new O:P) out of the
# way first, then rename/insert the new G:N analysis columns.
# ------------------------------------------------------------------

# old K column (meta-notes) -> new column O
$ws.Range("O1").Value = $ws.Range("K1").Value2
# old L column (meta-source) -> new column P
$ws.Range("P1").Value = $ws.Range("L1").Value2
$ws.Range("P3").Value = "ecoinvent 2.2"
$ws.Range("P4").Value = "ecoinvent 2.2"
$ws.Range("P9").Value = $ws.Range("L9").Value2
$ws.Range("P10").Value = $ws.Range("L10").Value2
$ws.Range("P11").Value = $ws.Range("L11").Value2
$ws.Range("P12").Value = $ws.Range("L12").Value2

# Clear the old L-column meta-source values (now duplicated into P)
$ws.Range("L9").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("L12").ClearContents()

# ------------------------------------------------------------------
# Row 1 headers: rename C/H/S/Ash "content" -> "%", and fill the new
# upstream-CO2 / biomass columns (K:N)
# ------------------------------------------------------------------
$ws.Range("G1").Value = "C %"
$ws.Range("H1").Value = "H %"
$ws.Range("I1").Value = "S %"
$ws.Range("J1").Value = "Ash %"
$ws.Range("K1").Value = "upstream CO2"
$ws.Range("L1").Value = "fresh biomass ratio"
$ws.Range("M1").Value = "fresh biomass C content"
$ws.Range("N1").Value = "biomass CO2 absorption"

# Row 2 meta-units for the new columns
$ws.Range("K2").Value = "t CO2 / t fuel"
$ws.Range("L2").Value = "t fresh / t fuel"
$ws.Range("M2").Value = "t c / t biomass"
$ws.Range("N2").Value = "t CO2 / t fresh"

# Row 3 (coal): upstream CO2 value + ecoinvent source (already set above)
$ws.Range("K3").Value = 0.064

# Row 4 (charcoal): upstream CO2 + fresh-biomass data + CO2-absorption formula
$ws.Range("K4").Value = 2.9
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 0.25
$ws.Range("N4").Formula = "=(44/12)*M4"
$ws.Range("N4").NumberFormat = "0.00"

# Row 8 (steam): LHV now mirrors HHV
$ws.Range("C8").Value = 2.77

# Row 12 (PCI coal): label fix -- was pointing at wrong shared string ("C content")
$ws.Range("A12").Value = "PCI coal"

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.125
$ws.Columns.Item(7).ColumnWidth = 6.375
$ws.Columns.Item(8).ColumnWidth = 6.25
$ws.Columns.Item(9).ColumnWidth = 6
$ws.Columns.Item(10).ColumnWidth = 8
$ws.Columns.Item(11).ColumnWidth = 11.25
$ws.Columns.Item(12).ColumnWidth = 18.625
$ws.Columns.Item(13).ColumnWidth = 12
$ws.Columns.Item(14).ColumnWidth = 13.5

# ------------------------------------------------------------------
# Header row formatting: wrap text + taller row
# ------------------------------------------------------------------
$ws.Range("A1:P1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 35.25

# ------------------------------------------------------------------
# Freeze the header row + first column, scrolled so K:P are in view
# ------------------------------------------------------------------
$ws.Range("K2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("O20").Select()
